$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 639  # was 1800
$ws.Range("J58").Value = 962  # was 3284
$ws.Range("L58").Value = 2886  # was 9852
$ws.Range("N58").Value = -3186  # was -10152
$ws.Range("H74").Value = 0  # was 2284.2856
$ws.Range("I74").Value = 0  # was 2000
$ws.Range("J74").Value = 0  # was 3990
$ws.Range("K74").Value = 0  # was 2000
$ws.Range("L74").Value = 0  # was 3990
$ws.Range("M74").ClearContents()  # was -1064
$ws.Range("N74").ClearContents()  # was -5862
$ws.Range("H77").Value = 0  # was 2284.2856
$ws.Range("I77").Value = 0  # was 2000
$ws.Range("J77").Value = 0  # was 3990
$ws.Range("K77").Value = 0  # was 10000
$ws.Range("L77").Value = 0  # was 19950
$ws.Range("M77").ClearContents()  # was -5320
$ws.Range("N77").ClearContents()  # was -29310
$ws.Range("H88").Value = 3332  # was 3998.3333
$ws.Range("J88").Value = 3332  # was 3998.3333
$ws.Range("L88").Value = 3332  # was 3998.3333
$ws.Range("N88").Value = -4144  # was -4810.3333
$ws.Range("H91").Value = 3332  # was 3998.3333
$ws.Range("J91").Value = 3332  # was 3998.3333
$ws.Range("L91").Value = 3332  # was 3998.3333
$ws.Range("N91").Value = -6140  # was -6806.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700  # was 925
$ws.Range("I2").Value = 300  # was 350
$ws.Range("K2").Value = 300  # was 350
$ws.Range("M2").Value = -187  # was -237
$ws.Range("H24").Value = 35000  # was 5023695
$ws.Range("J24").Value = 35000  # was 5023695
$ws.Range("L24").Value = 35000  # was 5023695
$ws.Range("N24").Value = -35748  # was -5024443
$ws.Range("H32").Value = 4229.2  # was 6571.5
$ws.Range("I32").Value = 3921.3333  # was 6485.8
$ws.Range("K32").Value = 3921.3333  # was 6485.8
$ws.Range("M32").Value = -3634.3333  # was -6198.8
$ws.Range("H61").Value = 2090.125  # was 2217.2856
$ws.Range("I61").Value = 2090.125  # was 2217.2856
$ws.Range("K61").Value = 2090.125  # was 2217.2856
$ws.Range("M61").Value = -1878.125  # was -2005.2856
$ws.Range("H62").Value = 0  # was 50000
$ws.Range("J62").Value = 0  # was 50000
$ws.Range("L62").Value = 0  # was 50000
$ws.Range("N62").ClearContents()  # was -51248
$ws.Range("H65").Value = 0  # was 50000
$ws.Range("J65").Value = 0  # was 50000
$ws.Range("L65").Value = 0  # was 150000
$ws.Range("N65").ClearContents()  # was -156240
$ws.Range("H88").Value = 1745.8334  # was 1695.2
$ws.Range("J88").Value = 2242.25  # was 2323.3333
$ws.Range("L88").Value = 2242.25  # was 2323.3333
$ws.Range("N88").Value = -3054.25  # was -3135.3333
$ws.Range("H91").Value = 1745.8334  # was 1695.2
$ws.Range("J91").Value = 2242.25  # was 2323.3333
$ws.Range("L91").Value = 2242.25  # was 2323.3333
$ws.Range("N91").Value = -5050.25  # was -5131.3333
$ws.Range("H100").Value = 35000  # was 5023695
$ws.Range("J100").Value = 35000  # was 5023695
$ws.Range("L100").Value = 35000  # was 5023695
$ws.Range("N100").Value = -37164  # was -5025859
$ws.Range("H102").Value = 2299.8  # was 2699.8
$ws.Range("I102").Value = 1624.75  # was 1833
$ws.Range("J102").Value = 5000  # was 4000
$ws.Range("K102").Value = 1624.75  # was 1833
$ws.Range("L102").Value = 5000  # was 4000
$ws.Range("M102").Value = -2.75  # was -211
$ws.Range("N102").Value = -8244  # was -7244
$ws.Range("H116").Value = 700  # was 925
$ws.Range("I116").Value = 300  # was 350
$ws.Range("K116").Value = 300  # was 350
$ws.Range("M116").Value = 1994  # was 1944
$ws.Range("H132").Value = 805  # was 1056.2
$ws.Range("I132").Value = 805  # was 1056.2
$ws.Range("K132").Value = 2415  # was 3168.6
$ws.Range("M132").Value = 115  # was -638.6000000000004
$ws.Range("H136").Value = 2090.125  # was 2217.2856
$ws.Range("I136").Value = 2090.125  # was 2217.2856
$ws.Range("K136").Value = 6270.375  # was 6651.8568
$ws.Range("M136").Value = -3720.375  # was -4101.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700  # was 925
$ws.Range("I3").Value = 300  # was 350
$ws.Range("K3").Value = 300  # was 350
$ws.Range("M3").Value = -186  # was -236
$ws.Range("H6").Value = 0  # was 50000
$ws.Range("J6").Value = 0  # was 50000
$ws.Range("L6").Value = 0  # was 50000
$ws.Range("N6").ClearContents()  # was -50226
$ws.Range("H26").Value = 17500  # was 20000
$ws.Range("I26").Value = 17500  # was 20000
$ws.Range("K26").Value = 17500  # was 20000
$ws.Range("M26").Value = -17208  # was -19708
$ws.Range("H95").Value = 24109  # was 27767.8
$ws.Range("J95").Value = 24109  # was 27767.8
$ws.Range("L95").Value = 24109  # was 27767.8
$ws.Range("N95").Value = -29601  # was -33259.8
$ws.Range("H96").Value = 17833.334  # was 20000
$ws.Range("I96").Value = 17833.334  # was 20000
$ws.Range("K96").Value = 17833.334  # was 20000
$ws.Range("M96").Value = -15087.334  # was -17254
$ws.Range("H135").Value = 59450  # was 99000
$ws.Range("J135").Value = 59450  # was 99000
$ws.Range("L135").Value = 59450  # was 99000
$ws.Range("N135").Value = -69590  # was -109140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2758  # was 4382.6
$ws.Range("I16").Value = 2947.5  # was 3633.3333
$ws.Range("J16").Value = 2000  # was 5506.5
$ws.Range("K16").Value = 2947.5  # was 3633.3333
$ws.Range("L16").Value = 2000  # was 5506.5
$ws.Range("M16").Value = -2660.5  # was -3346.3333
$ws.Range("N16").Value = -2574  # was -6080.5
$ws.Range("H25").Value = 49  # was 475.0909
$ws.Range("I25").Value = 49  # was 500
$ws.Range("J25").Value = 0  # was 408.66666
$ws.Range("K25").Value = 49  # was 500
$ws.Range("L25").Value = 0  # was 408.66666
$ws.Range("M25").Value = 125  # was -326
$ws.Range("N25").ClearContents()  # was -756.66666
$ws.Range("H43").Value = 7825  # was 7991.6665
$ws.Range("J43").Value = 7825  # was 7991.6665
$ws.Range("L43").Value = 7825  # was 7991.6665
$ws.Range("N43").Value = -8193  # was -8359.666499999999
$ws.Range("H92").Value = 19999.5  # was 20000
$ws.Range("J92").Value = 19999.5  # was 20000
$ws.Range("L92").Value = 19999.5  # was 20000
$ws.Range("N92").Value = -24991.5  # was -24992
$ws.Range("H101").Value = 7825  # was 7991.6665
$ws.Range("J101").Value = 7825  # was 7991.6665
$ws.Range("L101").Value = 7825  # was 7991.6665
$ws.Range("N101").Value = -14315  # was -14481.6665
$ws.Range("H113").Value = 2758  # was 4382.6
$ws.Range("I113").Value = 2947.5  # was 3633.3333
$ws.Range("J113").Value = 2000  # was 5506.5
$ws.Range("K113").Value = 2947.5  # was 3633.3333
$ws.Range("L113").Value = 2000  # was 5506.5
$ws.Range("M113").Value = -777.5  # was -1463.3333
$ws.Range("N113").Value = -6340  # was -9846.5
$ws.Range("H132").Value = 4871.75  # was 4497.1113
$ws.Range("I132").Value = 4192  # was 3743.3333
$ws.Range("K132").Value = 12576  # was 11229.9999
$ws.Range("M132").Value = -10046  # was -8699.999899999999
$ws.Range("H134").Value = 4975.0835  # was 5390.1
$ws.Range("I134").Value = 1140.4  # was 1225.5
$ws.Range("J134").Value = 7714.143  # was 8166.5
$ws.Range("K134").Value = 3421.2  # was 3676.5
$ws.Range("L134").Value = 23142.429  # was 24499.5
$ws.Range("M134").Value = -886.2000000000003  # was -1141.5
$ws.Range("N134").Value = -28212.429  # was -29569.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 724.1  # was 1371.3125
$ws.Range("J34").Value = 934.4  # was 1761.091
$ws.Range("L34").Value = 2803.2  # was 5283.272999999999
$ws.Range("N34").Value = -2971.2  # was -5451.272999999999
$ws.Range("H39").Value = 500  # was 450
$ws.Range("I39").Value = 500  # was 450
$ws.Range("J39").Value = 500  # was 0
$ws.Range("K39").Value = 1500  # was 1350
$ws.Range("L39").Value = 1500  # was 0
$ws.Range("M39").Value = -1206  # was -1056
$ws.Range("N39").Value = -2088  # was None
$ws.Range("H55").Value = 126512.375  # was 113844.336
$ws.Range("J55").Value = 2424.75  # was 4439.8
$ws.Range("L55").Value = 7274.25  # was 13319.4
$ws.Range("N55").Value = -7628.25  # was -13673.4
$ws.Range("H131").Value = 1654.8334  # was 1704.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 9566.429  # was 11081.667
$ws.Range("I9").Value = 1988.3334  # was 2745
$ws.Range("K9").Value = 1988.3334  # was 2745
$ws.Range("M9").Value = -1818.3334  # was -2575
$ws.Range("H102").Value = 1956.125  # was 2806.3333
$ws.Range("I102").Value = 1956.125  # was 2806.3333
$ws.Range("K102").Value = 1956.125  # was 2806.3333
$ws.Range("M102").Value = -334.125  # was -1184.3333
$ws.Range("H132").Value = 1819.1666  # was 1802.5
$ws.Range("I132").Value = 1913  # was 1784.5454
$ws.Range("J132").Value = 1350  # was 2000
$ws.Range("K132").Value = 5739  # was 5353.6362
$ws.Range("L132").Value = 4050  # was 6000
$ws.Range("M132").Value = -3209  # was -2823.6362
$ws.Range("N132").Value = -9110  # was -11060
$ws.Range("H133").Value = 100000  # was 0
$ws.Range("J133").Value = 100000  # was 0
$ws.Range("L133").Value = 100000  # was 0
$ws.Range("N133").Value = -110120  # was None

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 6738.2856  # was 6521
$ws.Range("I30").Value = 191.66667  # was 878.5714
$ws.Range("K30").Value = 191.66667  # was 878.5714
$ws.Range("M30").Value = -83.66667000000001  # was -770.5714
$ws.Range("H35").Value = 2528.5  # was 1903
$ws.Range("I35").Value = 2026.3334  # was 1370
$ws.Range("K35").Value = 2026.3334  # was 1370
$ws.Range("M35").Value = -1690.3334  # was -1034
$ws.Range("H68").Value = 2848.3333  # was 3850
$ws.Range("J68").Value = 847.5  # was 850
$ws.Range("L68").Value = 847.5  # was 850
$ws.Range("N68").Value = -2345.5  # was -2348
$ws.Range("H71").Value = 2848.3333  # was 3850
$ws.Range("J71").Value = 847.5  # was 850
$ws.Range("L71").Value = 4237.5  # was 4250
$ws.Range("N71").Value = -11725.5  # was -11738
$ws.Range("H76").Value = 23166.666  # was 23499.75
$ws.Range("J76").Value = 24500  # was 24499.5
$ws.Range("L76").Value = 24500  # was 24499.5
$ws.Range("N76").Value = -25176  # was -25175.5
$ws.Range("H79").Value = 23166.666  # was 23499.75
$ws.Range("J79").Value = 24500  # was 24499.5
$ws.Range("L79").Value = 24500  # was 24499.5
$ws.Range("N79").Value = -26840  # was -26839.5
$ws.Range("H122").Value = 4970.769  # was 4956.933
$ws.Range("I122").Value = 4692.8184  # was 4719.615
$ws.Range("K122").Value = 14078.4552  # was 14158.845
$ws.Range("M122").Value = -11628.4552  # was -11708.845
$ws.Range("H127").Value = 35607.25  # was 37476.332
$ws.Range("J127").Value = 35607.25  # was 37476.332
$ws.Range("L127").Value = 35607.25  # was 37476.332
$ws.Range("N127").Value = -45527.25  # was -47396.332
$ws.Range("H132").Value = 1000  # was 1500
$ws.Range("I132").Value = 1000  # was 1500
$ws.Range("K132").Value = 3000  # was 4500
$ws.Range("M132").Value = -470  # was -1970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 39271  # was 0
$ws.Range("J68").Value = 39271  # was 0
$ws.Range("L68").Value = 39271  # was 0
$ws.Range("N68").Value = -40893  # was None
$ws.Range("H71").Value = 39271  # was 0
$ws.Range("J71").Value = 39271  # was 0
$ws.Range("L71").Value = 117813  # was 0
$ws.Range("N71").Value = -125925  # was None
$ws.Range("H75").Value = 35000  # was 0
$ws.Range("J75").Value = 35000  # was 0
$ws.Range("L75").Value = 35000  # was 0
$ws.Range("N75").Value = -36872  # was None
$ws.Range("H78").Value = 35000  # was 0
$ws.Range("J78").Value = 35000  # was 0
$ws.Range("L78").Value = 105000  # was 0
$ws.Range("N78").Value = -114360  # was None
$ws.Range("H94").Value = 29665  # was 30000
$ws.Range("J94").Value = 29665  # was 30000
$ws.Range("L94").Value = 29665  # was 30000
$ws.Range("N94").Value = -31467  # was -31802
$ws.Range("H132").Value = 1659  # was 1371.3572
$ws.Range("I132").Value = 1273.75  # was 1066.5834
$ws.Range("K132").Value = 3821.25  # was 3199.7502
$ws.Range("M132").Value = -1291.25  # was -669.7501999999999
